$d = $word.ActiveDocument

# The sentence currently reads (across several runs, with stray grammar-check
# markers from a prior proofing pass):
#   "... product capabilities and" [" platform"][" "][gramStart]["integration"][gramEnd][" best practices."]
# The edit folds " " + "integration" + " best practices." into a single run
# and drops the now-stale <w:proofErr .../> gramStart/gramEnd pair, while
# leaving the " platform" run itself untouched.

# Locate the end of the " platform" run - that's where we'll temporarily work.
$plat = $d.Content
$plat.Find.Execute(" platform", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $plat.End

# Insert a one-character bold "barrier" right after " platform". Re-writing text
# that spans the proofErr markers makes the engine coalesce the touched runs
# together; without a formatting break immediately to the left, that coalescing
# would also swallow the preceding " platform" run. The bold marker prevents
# that by making the left neighbour's formatting visibly different during the edit.
$marker = $d.Range($insertPos, $insertPos)
$marker.InsertBefore("Z")
$markerRange = $d.Range($insertPos, $insertPos + 1)
$markerRange.Font.Bold = $true

# Re-apply the same text over itself across the old proofErr span. This is a
# no-op textually, but it forces Word to rebuild the runs it touches, which
# clears out the obsolete gramStart/gramEnd proofErr markers and merges
# "Z" + "integration" + " best practices." into one run (blocked from merging
# further left by the bold marker).
$rng = $d.Content
$rng.Find.Execute("Z integration best practices.", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "Z integration best practices.", 2)

# Remove the temporary marker character again.
$zRange = $d.Range($insertPos, $insertPos + 1)
$zRange.Text = ""

# Clear the bold formatting the marker carried so the merged run ends up as
# plain text, matching the rest of the sentence.
$finalRange = $d.Content
$finalRange.Find.Execute(" integration best practices.", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
$finalRange.Font.Bold = $false
